# cucc.xlsx edit: add a "Munka1" worksheet and extend the frequency table
# on Sheet1 with 3 more repetitions of the plate-number block (rows 71-97),
# matching the existing pattern used for rows 26-70.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Make sure Sheet1 is the active sheet before we start (matches the
# original file where Sheet1's view has tabSelected).
$sheet1.Activate()

# --- 1. Add the new empty "Munka1" worksheet, placed after Sheet1 -------
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Munka1"

# Keep Sheet1 as the active/selected sheet after inserting the new one.
$sheet1.Activate()

# --- 2. Extend Sheet1's data with 3 more cycles of the plate pattern ----
# The existing table (rows 26-70) cycles through 9 plate numbers:
# ABC-120, ABC-121, ABC-122, ABC-123, ABC-124, ABC-125, ABC-126, ABC-127, ABC-128
# Rows 71-97 continue that exact same 9-row cycle.
$plates = @("ABC-120", "ABC-121", "ABC-122", "ABC-123", "ABC-124", "ABC-125", "ABC-126", "ABC-127", "ABC-128")

for ($row = 71; $row -le 97; $row++) {
    $plate = $plates[($row - 71) % 9]
    $sheet1.Range("A$row").Value = $plate
    $sheet1.Range("B$row").Formula = "=(COUNTIF(A:A,A$row)-COUNTIF(A2,A$row))"
}
